$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new checklist row (row 24) at the end of the list, reusing the same
# row formatting as the row above it (row 23) by copying the whole row and
# inserting it right after, then overwriting its text.
$ws.Rows("23:23").Copy()
$ws.Rows("24:24").Insert()

$ws.Range("A24").Value = "fel medd när man försöker reg  användare 2 ggr"

# Move/restore the active selection like in the authored change.
$ws.Range("A28").Select()
